$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Update row 2 values: Cedula, Nombre, Telefono
$ws.Range("A2").Value = 9332945823
$ws.Range("B2").Value = "pedro picapiedra"
$ws.Range("C2").Value = 22691392

# Remove row 3 entirely (was Daniel / 103491814 / 3052076540)
$ws.Rows.Item(3).Delete()

# Adjust column widths to match target layout
# (ColumnWidth uses Excel's internal measurement which is offset from the
# raw OOXML <col width=".."> value by ~0.833; subtract that so the saved
# XML lands on exactly 12 / 18 / 10.)
$ws.Columns.Item(1).ColumnWidth = 11.166666667
$ws.Columns.Item(2).ColumnWidth = 17.166666667
$ws.Columns.Item(3).ColumnWidth = 9.166666667
